$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Influenza")
$ws.Select()

$rng = $ws.Range("I2:I6")
$rng.NumberFormat = "@"

$ws.Range("I2").Value = "01-Dec-2024"
$ws.Range("I3").Value = "31-Jan-2025"
$ws.Range("I4").Value = "01-Dec-2024"
$ws.Range("I5").Value = "31-Jan-2025"
$ws.Range("I6").Value = "01-Dec-2024"
